# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff event:
#   - file "4858cf89-ad94-41ff-be6a-4a1687b2b73d.md" moves from
#     "In Translation" to "Ready for handoff"
#   - the Latest Handoff timestamps for that file, plus
#     "0e809059-bc82-495f-9611-a1b7a42af53c.md" and
#     "d489ef34-099e-4119-97fe-6596e3553a26.md" are refreshed
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn = $wb.Worksheets.Item(2)
$dede = $wb.Worksheets.Item(3)

# --- Overview sheet ---
# Row 6: 0e809059-bc82-495f-9611-a1b7a42af53c.md
$overview.Range("D6").Value = "2016-23-20 16:23:57"

# Row 9: 4858cf89-ad94-41ff-be6a-4a1687b2b73d.md -> Ready for handoff
$overview.Range("B9").Value = "Ready for handoff"
$overview.Range("C9").Value = "Ready for handoff"
$overview.Range("D9").Value = "2016-23-20 16:23:57"

# Row 10: d489ef34-099e-4119-97fe-6596e3553a26.md
$overview.Range("D10").Value = "2016-23-20 16:23:57"

# --- zh-cn sheet ---
$zhcn.Range("E6").Value = "2016-03-20 16:23:54"

$zhcn.Range("C9").Value = "Ready for handoff"
$zhcn.Range("E9").Value = "2016-03-20 16:23:54"

$zhcn.Range("E10").Value = "2016-03-20 16:23:54"

# --- de-de sheet ---
$dede.Range("E6").Value = "2016-03-20 16:23:57"

$dede.Range("C9").Value = "Ready for handoff"
$dede.Range("E9").Value = "2016-03-20 16:23:57"

$dede.Range("E10").Value = "2016-03-20 16:23:57"
